$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 4307.5
$ws.Range("I48").Value = 4000
$ws.Range("K48").Value = 12000
$ws.Range("M48").Value = -11708
$ws.Range("H56").Value = 4307.5
$ws.Range("I56").Value = 4000
$ws.Range("K56").Value = 12000
$ws.Range("M56").Value = -11466
$ws.Range("H62").Value = 2687
$ws.Range("I62").Value = 2687
$ws.Range("K62").Value = 2687
$ws.Range("M62").Value = -2063
$ws.Range("H65").Value = 2687
$ws.Range("I65").Value = 2687
$ws.Range("K65").Value = 13435
$ws.Range("M65").Value = -10315
$ws.Range("H99").Value = 1197.25
$ws.Range("I99").Value = 307.6
$ws.Range("K99").Value = 922.8000000000001
$ws.Range("M99").Value = 575.1999999999999
$ws.Range("H135").Value = 3003.8
$ws.Range("I135").Value = 3286.125
$ws.Range("J135").Value = 1874.5
$ws.Range("K135").Value = 29575.125
$ws.Range("L135").Value = 16870.5
$ws.Range("M135").Value = -27040.125
$ws.Range("N135").Value = -21940.5
$ws.Range("H137").Value = 1645.4
$ws.Range("I137").Value = 1420.4783
$ws.Range("J137").Value = 2384.4285
$ws.Range("K137").Value = 4261.4349
$ws.Range("L137").Value = 7153.2855
$ws.Range("M137").Value = -1711.4349
$ws.Range("N137").Value = -12253.2855

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7296.125
$ws.Range("I61").Value = 2829.1765
$ws.Range("J61").Value = 18144.428
$ws.Range("K61").Value = 2829.1765
$ws.Range("L61").Value = 18144.428
$ws.Range("M61").Value = -2617.1765
$ws.Range("N61").Value = -18568.428
$ws.Range("H74").Value = 2196.2368
$ws.Range("I74").Value = 2283.125
$ws.Range("K74").Value = 2283.125
$ws.Range("M74").Value = -1409.125
$ws.Range("H77").Value = 2196.2368
$ws.Range("I77").Value = 2283.125
$ws.Range("K77").Value = 11415.625
$ws.Range("M77").Value = -7047.625
$ws.Range("H96").Value = 45000
$ws.Range("J96").Value = 45000
$ws.Range("L96").Value = 45000
$ws.Range("N96").Value = -50492
$ws.Range("H97").Value = 1473.3889
$ws.Range("I97").Value = 908.7692
$ws.Range("J97").Value = 2941.4
$ws.Range("K97").Value = 908.7692
$ws.Range("L97").Value = 2941.4
$ws.Range("M97").Value = -412.7692
$ws.Range("N97").Value = -3933.4
$ws.Range("H102").Value = 2961.0833
$ws.Range("I102").Value = 3057.5454
$ws.Range("J102").Value = 1900
$ws.Range("K102").Value = 3057.5454
$ws.Range("L102").Value = 1900
$ws.Range("M102").Value = -1435.5454
$ws.Range("N102").Value = -5144
$ws.Range("H122").Value = 2651.3333
$ws.Range("I122").Value = 1460.0526
$ws.Range("J122").Value = 4709
$ws.Range("K122").Value = 4380.1578
$ws.Range("L122").Value = 14127
$ws.Range("M122").Value = -1930.1578
$ws.Range("N122").Value = -19027
$ws.Range("H132").Value = 2030.619
$ws.Range("I132").Value = 1837.4375
$ws.Range("J132").Value = 2648.8
$ws.Range("K132").Value = 5512.3125
$ws.Range("L132").Value = 7946.400000000001
$ws.Range("M132").Value = -2982.3125
$ws.Range("N132").Value = -13006.4
$ws.Range("H136").Value = 7296.125
$ws.Range("I136").Value = 2829.1765
$ws.Range("J136").Value = 18144.428
$ws.Range("K136").Value = 8487.529500000001
$ws.Range("L136").Value = 54433.284
$ws.Range("M136").Value = -5937.529500000001
$ws.Range("N136").Value = -59533.284

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 60000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 60000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 60000
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -60630
$ws.Range("H79").Value = 60000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 60000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 60000
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -62184
$ws.Range("H94").Value = 18847.875
$ws.Range("I94").Value = 10021.375
$ws.Range("J94").Value = 27674.375
$ws.Range("K94").Value = 10021.375
$ws.Range("L94").Value = 27674.375
$ws.Range("M94").Value = -9570.375
$ws.Range("N94").Value = -28576.375
$ws.Range("H105").Value = 3305.923
$ws.Range("J105").Value = 1932.3334
$ws.Range("L105").Value = 1932.3334
$ws.Range("N105").Value = -5426.3334
$ws.Range("H107").Value = 18524.625
$ws.Range("I107").Value = 20442.428
$ws.Range("K107").Value = 20442.428
$ws.Range("M107").Value = -18522.428

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 159.75
$ws.Range("I7").Value = 189.66667
$ws.Range("J7").Value = 129.83333
$ws.Range("K7").Value = 189.66667
$ws.Range("L7").Value = 129.83333
$ws.Range("M7").Value = -76.66667000000001
$ws.Range("N7").Value = -355.83333
$ws.Range("H16").Value = 481.25
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 2730.111
$ws.Range("I31").Value = 1996.7142
$ws.Range("J31").Value = 3196.818
$ws.Range("K31").Value = 1996.7142
$ws.Range("L31").Value = 3196.818
$ws.Range("M31").Value = -1701.7142
$ws.Range("N31").Value = -3786.818
$ws.Range("H34").Value = 2730.111
$ws.Range("I34").Value = 1996.7142
$ws.Range("J34").Value = 3196.818
$ws.Range("K34").Value = 1996.7142
$ws.Range("L34").Value = 3196.818
$ws.Range("M34").Value = -1794.7142
$ws.Range("N34").Value = -3600.818
$ws.Range("H99").Value = 2049.5715
$ws.Range("I99").Value = 2100
$ws.Range("K99").Value = 2100
$ws.Range("M99").Value = -602
$ws.Range("H107").Value = 455.66666
$ws.Range("I107").Value = 508.29413
$ws.Range("J107").Value = 327.85715
$ws.Range("K107").Value = 508.29413
$ws.Range("L107").Value = 327.85715
$ws.Range("M107").Value = 1411.70587
$ws.Range("N107").Value = -4167.85715
$ws.Range("H113").Value = 481.25
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 2049.5715
$ws.Range("I126").Value = 2100
$ws.Range("K126").Value = 6300
$ws.Range("M126").Value = -3830

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 383.21054
$ws.Range("I5").Value = 355.8125
$ws.Range("J5").Value = 529.3333
$ws.Range("K5").Value = 1067.4375
$ws.Range("L5").Value = 1587.9999
$ws.Range("M5").Value = -955.4375
$ws.Range("N5").Value = -1811.9999
$ws.Range("H64").Value = 10999
$ws.Range("J64").Value = 10999
$ws.Range("L64").Value = 32997
$ws.Range("N64").Value = -33537
$ws.Range("H67").Value = 10999
$ws.Range("J67").Value = 10999
$ws.Range("L67").Value = 32997
$ws.Range("N67").Value = -34869
$ws.Range("H75").Value = 910.3333
$ws.Range("I75").Value = 1125
$ws.Range("J75").Value = 803
$ws.Range("K75").Value = 3375
$ws.Range("L75").Value = 2409
$ws.Range("M75").Value = -2377
$ws.Range("N75").Value = -4405
$ws.Range("H78").Value = 910.3333
$ws.Range("I78").Value = 1125
$ws.Range("J78").Value = 803
$ws.Range("K78").Value = 10125
$ws.Range("L78").Value = 7227
$ws.Range("M78").Value = -5133
$ws.Range("N78").Value = -17211
$ws.Range("H107").Value = 614.4286
$ws.Range("J107").Value = 1171.8334
$ws.Range("L107").Value = 3515.5002
$ws.Range("N107").Value = -7355.5002
$ws.Range("H116").Value = 2964.5
$ws.Range("I116").Value = 2964.5
$ws.Range("K116").Value = 8893.5
$ws.Range("M116").Value = -5451.5
$ws.Range("H117").Value = 3366.75
$ws.Range("I117").Value = 490
$ws.Range("J117").Value = 4325.6665
$ws.Range("K117").Value = 1470
$ws.Range("L117").Value = 12976.9995
$ws.Range("M117").Value = 1972
$ws.Range("N117").Value = -19860.9995
$ws.Range("H135").Value = 383.21054
$ws.Range("I135").Value = 355.8125
$ws.Range("J135").Value = 529.3333
$ws.Range("K135").Value = 3202.3125
$ws.Range("L135").Value = 4763.9997
$ws.Range("M135").Value = -667.3125
$ws.Range("N135").Value = -9833.9997
$ws.Range("H137").Value = 7284.875
$ws.Range("I137").Value = 2812.3
$ws.Range("J137").Value = 14739.167
$ws.Range("K137").Value = 8436.900000000001
$ws.Range("L137").Value = 44217.501
$ws.Range("M137").Value = -3336.900000000001
$ws.Range("N137").Value = -54417.501

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7163.115
$ws.Range("I70").Value = 6943.316
$ws.Range("J70").Value = 7759.7144
$ws.Range("K70").Value = 6943.316
$ws.Range("L70").Value = 7759.7144
$ws.Range("M70").Value = -6673.316
$ws.Range("N70").Value = -8299.714400000001
$ws.Range("H73").Value = 7163.115
$ws.Range("I73").Value = 6943.316
$ws.Range("J73").Value = 7759.7144
$ws.Range("K73").Value = 6943.316
$ws.Range("L73").Value = 7759.7144
$ws.Range("M73").Value = -6007.316
$ws.Range("N73").Value = -9631.714400000001
$ws.Range("H97").Value = 1379.6072
$ws.Range("I97").Value = 429
$ws.Range("K97").Value = 429
$ws.Range("M97").Value = 67

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1549.3334
$ws.Range("I93").Value = 824.5
$ws.Range("J93").Value = 2999
$ws.Range("K93").Value = 824.5
$ws.Range("L93").Value = 2999
$ws.Range("M93").Value = 423.5
$ws.Range("N93").Value = -5495
$ws.Range("H132").Value = 2257.6333
$ws.Range("I132").Value = 2114.6667
$ws.Range("J132").Value = 2352.9443
$ws.Range("K132").Value = 6344.000100000001
$ws.Range("L132").Value = 7058.8329
$ws.Range("M132").Value = -3814.000100000001
$ws.Range("N132").Value = -12118.8329

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5635.8335
$ws.Range("I62").Value = 5205.7144
$ws.Range("J62").Value = 6238
$ws.Range("K62").Value = 5205.7144
$ws.Range("L62").Value = 6238
$ws.Range("M62").Value = -4581.7144
$ws.Range("N62").Value = -7486
$ws.Range("H65").Value = 5635.8335
$ws.Range("I65").Value = 5205.7144
$ws.Range("J65").Value = 6238
$ws.Range("K65").Value = 26028.572
$ws.Range("L65").Value = 31190
$ws.Range("M65").Value = -22908.572
$ws.Range("N65").Value = -37430
$ws.Range("H132").Value = 4957.3213
$ws.Range("I132").Value = 4774.5654
$ws.Range("K132").Value = 14323.6962
$ws.Range("M132").Value = -11793.6962

